$d = $word.ActiveDocument

# --- Paragraph 1: update the hidden ID bookmark text, its indentation, and add a paragraph border ---
$p1 = $d.Paragraphs(1)

# Replace the ID placeholder text (without touching the trailing space run yet)
# MatchWildcards=$false so the literal "*" characters in the search string are matched as-is.
$d.Content.Find.Execute("**ID__AFFARS_5333_topic_7__ID**", $true, $false, $false, $false, $false, $false, 1, $false, "**ID__AFFARS_5333_170__ID**", 2) | Out-Null

# Remove the now-orphaned trailing space run left in the paragraph
$p1 = $d.Paragraphs(1)
$endRng = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
if ($endRng.Text -eq " ") {
    $endRng.Delete()
}

# Update paragraph formatting: left indent 120 -> 225 twips (6pt -> 11.25pt)
$p1 = $d.Paragraphs(1)
$pf = $p1.Range.ParagraphFormat
$pf.LeftIndent = 11.25

# Add a paragraph border (space=5 on all sides, no line drawn)
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
